# Update crypto price/volume figures for the Tue Mar 14 08:53:20 UTC 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.384.88"
$ws.Range("E2").Value = "  +8.95%  "
$ws.Range("D3").Value = "1.677.01"
$ws.Range("E3").Value = "  +4.74%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'0.9999"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "'305.86"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.3696"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").Value = "'0.3422"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "'47.54"
$ws.Range("E9").Value = "  +12.96%  "
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'6.149"
$ws.Range("E13").Value = "  +4.28%  "
$ws.Range("D14").Value = "'20.10"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "1.675.45"
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "'0.9998"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'0.06659"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "'80.58"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "'16.43"
$ws.Range("E21").Value = "  +2.94%  "
$ws.Range("D22").Value = "'6.098"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").Value = "'12.19"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("D24").Value = "24.351.92"
$ws.Range("E24").Value = "  +8.57%  "
$ws.Range("D25").Value = "'2.430"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").Value = "'2.649"
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("D27").Value = "'152.09"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "'19.46"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "1.862.61"
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("D30").Value = "'127.85"
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("D31").Value = "'6.281"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").Value = "'4.050"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  +2.81%  "
$ws.Range("D34").Value = "'0.08436"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").Value = "'1.685"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "'12.21"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").Value = "'0.06382"
$ws.Range("E37").Value = "  +5.46%  "
$ws.Range("D38").Value = "'5.314"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").Value = "'8.721"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "'0.02316"
$ws.Range("D41").Value = "'1.237"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "'0.2085"
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("D43").Value = "'0.6089"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("D45").Value = "'3.746"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").Value = "'12.95"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "'0.5878"
$ws.Range("E47").Value = "  +3.84%  "
$ws.Range("D48").Value = "'126.08"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "'2.015"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("D50").Value = "'0.07145"
$ws.Range("E50").Value = "  +5.12%  "
$ws.Range("D51").Value = "'75.63"
$ws.Range("E51").Value = "  +3.09%  "
